$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unmerge the original merged header cells before moving content ---
$ws.Range("B3:B4").UnMerge()
$ws.Range("C3:E3").UnMerge()
$ws.Range("B11:B12").UnMerge()
$ws.Range("C11:E11").UnMerge()

# --- 2. Move the whole B3:F18 block (both tables) to I3:M18, keeping styles/values ---
$ws.Range("B3:F18").Copy($ws.Range("I3"))

# --- 3. Clear the old location entirely (values + formatting) ---
$ws.Range("B3:F18").Clear()

# Remove leftover blank rows created by the block copy (rows with no data in the source)
$ws.Range("I9:M10").Clear()
$ws.Range("I17:M17").Clear()

# --- 4. Re-create the merges at the new location ---
$ws.Range("I3:I4").Merge()
$ws.Range("J3:L3").Merge()
$ws.Range("I11:I12").Merge()
$ws.Range("J11:L11").Merge()

# --- 5. Fill in table 2 (rows 13-16) input values and formulas ---
$ws.Range("J13").Value = 100
$ws.Range("L14").Value = 140
$ws.Range("J15").Value = 20
$ws.Range("K15").Value = 130
$ws.Range("L15").Value = 10

$ws.Range("M13").Formula = "=SUM(J13:L13)"
$ws.Range("M14").Formula = "=SUM(J14:L14)"
$ws.Range("M15").Formula = "=SUM(J15:L15)"

$ws.Range("J16").Formula = "=SUM(J13:J15)"
$ws.Range("K16").Formula = "=SUM(K13:K15)"
$ws.Range("L16").Formula = "=SUM(L13:L15)"

$ws.Range("J18").Formula = "=SUMPRODUCT(J13:L15,J5:L7)"

# Tidy up stray blank cells introduced by the rectangular block-copy
$ws.Range("M8").Clear()
$ws.Range("M16").Clear()
$ws.Range("K18:M18").Clear()

Write-Host "table2 done"

# --- 6. Build table 3 (rows 21-28) as a copy of table 2 (rows 11-18) ---
$ws.Range("I11:M18").Copy($ws.Range("I21"))

# Tidy up stray blank cells introduced by the rectangular block-copy
$ws.Range("I27:M27").Clear()
$ws.Range("K28:M28").Clear()

$ws.Range("I21:I22").Merge()
$ws.Range("J21:L21").Merge()

# Overwrite the input cells with table 3's own values
$ws.Range("J23").Value = 100
$ws.Range("K23").ClearContents()
$ws.Range("L23").ClearContents()

$ws.Range("J24").Value = 20
$ws.Range("K24").ClearContents()
$ws.Range("L24").Value = 120

$ws.Range("J25").ClearContents()
$ws.Range("K25").Value = 130
$ws.Range("L25").Value = 30

$ws.Range("M23").Formula = "=SUM(J23:L23)"
$ws.Range("M24").Formula = "=SUM(J24:L24)"
$ws.Range("M25").Formula = "=SUM(J25:L25)"

$ws.Range("J26").Formula = "=SUM(J23:J25)"
$ws.Range("K26").Formula = "=SUM(K23:K25)"
$ws.Range("L26").Formula = "=SUM(L23:L25)"
$ws.Range("M26").Clear()

$ws.Range("J28").Formula = "=SUMPRODUCT(J23:L25,J5:L7)"

Write-Host "table3 done"

# --- 7. Row 31: single-space label ---
$ws.Range("L31").Value = " "

Write-Host "row31 done"

# --- 8. View state: zoom + selection on the newly-added table ---
$ws.Application.ActiveWindow.Zoom = 150
$ws.Range("I11:M16").Select()

Write-Host "view done"

# --- 9. Column widths follow the data to their new home (I / J:M) ---
$ws.Columns("I:I").ColumnWidth = 11.498697916666666
$ws.Columns("J:M").ColumnWidth = 11.666666666666666

Write-Host "cols done"
